# Updated symbol list on Mon Dec 26 09:56:06 UTC 2022 with GitHub Actions
# Apply the price/volume-label refresh described by the commit's diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a cell whose displayed text must remain a plain text
# string (these "Price" column cells are numeric-looking but are stored
# as text in the workbook, so we force text formatting before/after the
# assignment to avoid Excel auto-converting them to real numbers).
function Set-TextValue($Worksheet, $Address, $NewText) {
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $NewText
    $cell.Style = "Normal"
}

# --- Price column (D) updates -------------------------------------------------
Set-TextValue $ws "D2"  "243.04"
Set-TextValue $ws "D3"  "23.06"
Set-TextValue $ws "D4"  "5.402"
Set-TextValue $ws "D5"  "0.05922"
Set-TextValue $ws "D6"  "3.448"
Set-TextValue $ws "D7"  "6.527"
Set-TextValue $ws "D8"  "0.8107"
Set-TextValue $ws "D9"  "0.9095"
Set-TextValue $ws "D10" "0.1410"
Set-TextValue $ws "D11" "0.07329"
Set-TextValue $ws "D12" "0.03264"
Set-TextValue $ws "D13" "0.03044"
Set-TextValue $ws "D14" "0.09353"
Set-TextValue $ws "D15" "3.853"
Set-TextValue $ws "D16" "0.001559"
Set-TextValue $ws "D17" "0.04675"
Set-TextValue $ws "D18" "0.0005942"
Set-TextValue $ws "D19" "0.006075"
Set-TextValue $ws "D20" "0.004978"
Set-TextValue $ws "D21" "0.0009811"
Set-TextValue $ws "D22" "0.00009404"
Set-TextValue $ws "D27" "0.0002901"
Set-TextValue $ws "D40" "0.03965"
Set-TextValue $ws "D41" "0.006197"
Set-TextValue $ws "D43" "0.003001"
Set-TextValue $ws "D44" "0.008187"
Set-TextValue $ws "D45" "0.00005244"
Set-TextValue $ws "D47" "0.7503"
Set-TextValue $ws "D48" "0.002245"

# --- Volume(1h) column (E) label updates --------------------------------------
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
